$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 98.888885
$ws.Range("I2").Value = 65
$ws.Range("J2").Value = 166.66667
$ws.Range("K2").Value = 65
$ws.Range("L2").Value = 166.66667
$ws.Range("M2").Value = 48
$ws.Range("N2").Value = -392.66667

# Row 4
$ws.Range("H4").Value = 316.27777
$ws.Range("I4").Value = 174.16667
$ws.Range("J4").Value = 600.5
$ws.Range("K4").Value = 174.16667
$ws.Range("L4").Value = 600.5
$ws.Range("M4").Value = -60.16667000000001
$ws.Range("N4").Value = -828.5

# Row 6
$ws.Range("H6").Value = 531.6667
$ws.Range("J6").Value = 1666.6666
$ws.Range("L6").Value = 4999.9998
$ws.Range("N6").Value = -5223.9998

# Row 86
$ws.Range("H86").Value = 4820.3335
$ws.Range("J86").Value = 5855.778
$ws.Range("L86").Value = 5855.778
$ws.Range("N86").Value = -8101.778

# Row 89
$ws.Range("H89").Value = 4820.3335
$ws.Range("J89").Value = 5855.778
$ws.Range("L89").Value = 29278.89
$ws.Range("N89").Value = -40510.89

# Row 106
$ws.Range("H106").Value = 1435435.4
$ws.Range("I106").Value = 2503262
$ws.Range("J106").Value = 11666.667
$ws.Range("K106").Value = 2503262
$ws.Range("L106").Value = 11666.667
$ws.Range("M106").Value = -2502631
$ws.Range("N106").Value = -12928.667

# Row 111
$ws.Range("H111").Value = 1525.8
$ws.Range("I111").Value = 1543
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 4629
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = -1562
$ws.Range("N111").Value = -10634

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7704.9775
$ws.Range("I32").Value = 3865.3242
$ws.Range("K32").Value = 3865.3242
$ws.Range("M32").Value = -3578.3242

# Row 41
$ws.Range("H41").Value = 2990.8333
$ws.Range("I41").Value = 2990.8333
$ws.Range("K41").Value = 2990.8333
$ws.Range("M41").Value = -2576.8333

# Row 74
$ws.Range("H74").Value = 41868.6
$ws.Range("I74").Value = 44590.87
$ws.Range("J74").Value = 10562.5
$ws.Range("K74").Value = 44590.87
$ws.Range("L74").Value = 10562.5
$ws.Range("M74").Value = -43716.87
$ws.Range("N74").Value = -12310.5

# Row 77
$ws.Range("H77").Value = 41868.6
$ws.Range("I77").Value = 44590.87
$ws.Range("J77").Value = 10562.5
$ws.Range("K77").Value = 222954.35
$ws.Range("L77").Value = 52812.5
$ws.Range("M77").Value = -218586.35
$ws.Range("N77").Value = -61548.5

$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Range("H10").Value = 2752.5
$ws.Range("I10").Value = 2752.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 2752.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -2612.5
$ws.Range("N10").ClearContents()

# Row 100
$ws.Range("H100").Value = 23824.834
$ws.Range("J100").Value = 23824.834
$ws.Range("L100").Value = 23824.834
$ws.Range("N100").Value = -25988.834

# Row 107
$ws.Range("H107").Value = 925.7143
$ws.Range("I107").Value = 810.8333
$ws.Range("J107").Value = 1615
$ws.Range("K107").Value = 810.8333
$ws.Range("L107").Value = 1615
$ws.Range("M107").Value = 1109.1667
$ws.Range("N107").Value = -5455

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1847.9193
$ws.Range("I31").Value = 930.5135
$ws.Range("J31").Value = 3205.68
$ws.Range("K31").Value = 930.5135
$ws.Range("L31").Value = 3205.68
$ws.Range("M31").Value = -635.5135
$ws.Range("N31").Value = -3795.68

# Row 34
$ws.Range("H34").Value = 1847.9193
$ws.Range("I34").Value = 930.5135
$ws.Range("J34").Value = 3205.68
$ws.Range("K34").Value = 930.5135
$ws.Range("L34").Value = 3205.68
$ws.Range("M34").Value = -728.5135
$ws.Range("N34").Value = -3609.68

# Row 96
$ws.Range("H96").Value = 19163.637
$ws.Range("J96").Value = 19163.637
$ws.Range("L96").Value = 19163.637
$ws.Range("N96").Value = -24655.637

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 476.2381
$ws.Range("I4").Value = 133.4
$ws.Range("J4").Value = 1333.3334
$ws.Range("K4").Value = 400.2
$ws.Range("L4").Value = 4000.0002
$ws.Range("M4").Value = -288.2
$ws.Range("N4").Value = -4224.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 151.11111
$ws.Range("I2").Value = 68.5
$ws.Range("J2").Value = 192.41667
$ws.Range("K2").Value = 68.5
$ws.Range("L2").Value = 192.41667
$ws.Range("M2").Value = 44.5
$ws.Range("N2").Value = -418.41667

# Row 4
$ws.Range("H4").Value = 3998
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# Row 9
$ws.Range("H9").Value = 20995
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 20995
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 20995
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -21335

# Row 35
$ws.Range("H35").Value = 3015
$ws.Range("I35").Value = 3015
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 3015
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -2717
$ws.Range("N35").ClearContents()

# Row 70
$ws.Range("H70").Value = 5114.0454
$ws.Range("I70").Value = 4944.4443
$ws.Range("J70").Value = 5877.25
$ws.Range("K70").Value = 4944.4443
$ws.Range("L70").Value = 5877.25
$ws.Range("M70").Value = -4674.4443
$ws.Range("N70").Value = -6417.25

# Row 73
$ws.Range("H73").Value = 5114.0454
$ws.Range("I73").Value = 4944.4443
$ws.Range("J73").Value = 5877.25
$ws.Range("K73").Value = 4944.4443
$ws.Range("L73").Value = 5877.25
$ws.Range("M73").Value = -4008.4443
$ws.Range("N73").Value = -7749.25

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 20800
$ws.Range("J2").Value = 3500
$ws.Range("L2").Value = 3500
$ws.Range("N2").Value = -3724

# Row 3
$ws.Range("H3").Value = 750
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 666.6667
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 666.6667
$ws.Range("M3").Value = -886
$ws.Range("N3").Value = -894.6667

# Row 4
$ws.Range("H4").Value = 90000
$ws.Range("I4").Value = 90000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 90000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -89887
$ws.Range("N4").ClearContents()

# Row 113
$ws.Range("H113").Value = 729.7222
$ws.Range("I113").Value = 576.8889
$ws.Range("J113").Value = 882.55554
$ws.Range("K113").Value = 1730.6667
$ws.Range("L113").Value = 1730.6667
$ws.Range("M113").Value = 439.3332999999998
$ws.Range("N113").Value = -6070.6667
